$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header cleanup: drop the "UpdateAction" / "UpdateComment" columns ---
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = ""

# --- Row 2: K/L (old "CM - Pricing" / "test pricing") no longer apply ---
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""

# --- Row 2: G/H now hold the "DIO - RMA" aging test case instead of "CM - Cost Structure" ---
$ws.Range("G2").Value = "DIO - RMA"
$ws.Range("H2").Value = "rma test"

# --- New row 3: the old "CM - Cost Structure" aging test case, moved down ---
$ws.Range("G3").Value = "CM - Cost Structure"
$ws.Range("H3").Value = "test cost structure"
$ws.Range("G2:H2").Copy() | Out-Null
$ws.Range("G3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# --- New row 4: the old "CM - Pricing" aging test case, moved down ---
$ws.Range("G4").Value = "CM - Pricing"
$ws.Range("H4").Value = "pricing"
$ws.Range("G2:H2").Copy() | Out-Null
$ws.Range("G4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# --- New row 6: Top 100 aging UnderperformingSKU testcase ---
$ws.Range("A6").Value = "iinventory_management_action_planning"
$ws.Range("B6").Value = "9YG957"
$ws.Range("C6").Value = "MLYU3AM/A"
$ws.Range("D6").Value = "ROBERT BROWN"
$ws.Range("E6").Value = "JUNIPER"
$ws.Range("F6").Value = "MD"
$ws.Range("G6").Value = "DIO - RMA"
$ws.Range("H6").Value = "rma test"

# A6 mirrors the style already used by A2 (bold-ish vertical-center default style)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# G6/H6 mirror the style already used on G2/H2
$ws.Range("G2:H2").Copy() | Out-Null
$ws.Range("G6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# C6 / E6 use a dedicated Roboto / #333333 font
$ws.Range("C6").Font.Color = 3355443
$ws.Range("C6").Font.Name = "Roboto"
$ws.Range("C6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- Selection moved to K17 ---
$ws.Range("K17").Select() | Out-Null
